$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 345.33334
$ws.Range("I58").Value = 248.55556
$ws.Range("J58").Value = 635.6667
$ws.Range("K58").Value = 745.66668
$ws.Range("L58").Value = 1907.0001
$ws.Range("M58").Value = -595.66668
$ws.Range("N58").Value = -2207.0001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1902.3226
$ws.Range("I98").Value = 1715.7037
$ws.Range("J98").Value = 3162
$ws.Range("K98").Value = 1715.7037
$ws.Range("L98").Value = 3162
$ws.Range("M98").Value = -217.7037
$ws.Range("N98").Value = -6158

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1902.3226
$ws.Range("I122").Value = 1715.7037
$ws.Range("J122").Value = 3162
$ws.Range("K122").Value = 5147.1111
$ws.Range("L122").Value = 9486
$ws.Range("M122").Value = -2697.1111
$ws.Range("N122").Value = -14386

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2845.2856
$ws.Range("I125").Value = 2994.4
$ws.Range("J125").Value = 2472.5
$ws.Range("K125").Value = 26949.6
$ws.Range("L125").Value = 22252.5
$ws.Range("M125").Value = -24489.6
$ws.Range("N125").Value = -27172.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1808.1428
$ws.Range("I129").Value = 1484.1818
$ws.Range("J129").Value = 2996
$ws.Range("K129").Value = 4452.5454
$ws.Range("L129").Value = 8988
$ws.Range("M129").Value = 547.4546
$ws.Range("N129").Value = -18988

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 6139.3335
$ws.Range("I131").Value = 1635.3334
$ws.Range("K131").Value = 4906.0002
$ws.Range("M131").Value = 133.9997999999996

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2959.8667
$ws.Range("I137").Value = 2323.111
$ws.Range("J137").Value = 3915
$ws.Range("K137").Value = 6969.333
$ws.Range("L137").Value = 11745
$ws.Range("M137").Value = -4419.333
$ws.Range("N137").Value = -16845

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2506.2903
$ws.Range("I138").Value = 1360.3334
$ws.Range("J138").Value = 3580.625
$ws.Range("K138").Value = 4081.0002
$ws.Range("L138").Value = 10741.875
$ws.Range("M138").Value = 1058.9998
$ws.Range("N138").Value = -21021.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 947.10345
$ws.Range("I141").Value = 983.1852
$ws.Range("J141").Value = 460
$ws.Range("K141").Value = 2949.5556
$ws.Range("L141").Value = 1380
$ws.Range("M141").Value = 2230.4444
$ws.Range("N141").Value = -11740

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3888.1875
$ws.Range("I32").Value = 3878.1428
$ws.Range("K32").Value = 3878.1428
$ws.Range("M32").Value = -3591.1428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 250002260
$ws.Range("I61").Value = 500001250
$ws.Range("K61").Value = 500001250
$ws.Range("M61").Value = -500001038

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4352137
$ws.Range("I132").Value = 5267505.5
$ws.Range("J132").Value = 4136.25
$ws.Range("K132").Value = 15802516.5
$ws.Range("L132").Value = 12408.75
$ws.Range("M132").Value = -15799986.5
$ws.Range("N132").Value = -17468.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 250002260
$ws.Range("I136").Value = 500001250
$ws.Range("K136").Value = 1500003750
$ws.Range("M136").Value = -1500001200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 39231656
$ws.Range("I134").Value = 39231656
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 117694968
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -117692433

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").ClearContents()
$ws.Range("N25").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 50011480
$ws.Range("I58").Value = 125025740
$ws.Range("J58").Value = 1970.6666
$ws.Range("K58").Value = 125025740
$ws.Range("L58").Value = 1970.6666
$ws.Range("M58").Value = -125025537
$ws.Range("N58").Value = -2376.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 75000
$ws.Range("I59").Value = 16667
$ws.Range("J59").Value = 133333
$ws.Range("K59").Value = 16667
$ws.Range("L59").Value = 133333
$ws.Range("M59").Value = -15522
$ws.Range("N59").Value = -135623

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3922.6667
$ws.Range("I99").Value = 3907.2
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 3907.2
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = -2409.2
$ws.Range("N99").Value = -6996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2603.1667
$ws.Range("I122").Value = 2603.1667
$ws.Range("K122").Value = 7809.500100000001
$ws.Range("M122").Value = -5359.500100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3922.6667
$ws.Range("I126").Value = 3907.2
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 11721.6
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -9251.599999999999
$ws.Range("N126").Value = -16940

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 17858862
$ws.Range("I132").Value = 20001766
$ws.Range("J132").Value = 1336
$ws.Range("K132").Value = 60005298
$ws.Range("L132").Value = 4008
$ws.Range("M132").Value = -60002768
$ws.Range("N132").Value = -9068

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 7145055.5
$ws.Range("I134").Value = 9261521
$ws.Range("J134").Value = 1986.125
$ws.Range("K134").Value = 27784563
$ws.Range("L134").Value = 5958.375
$ws.Range("M134").Value = -27782028
$ws.Range("N134").Value = -11028.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 50011480
$ws.Range("I136").Value = 125025740
$ws.Range("J136").Value = 1970.6666
$ws.Range("K136").Value = 375077220
$ws.Range("L136").Value = 5911.9998
$ws.Range("M136").Value = -375074670
$ws.Range("N136").Value = -11011.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 592.8333
$ws.Range("J92").Value = 716.75
$ws.Range("L92").Value = 2150.25
$ws.Range("N92").Value = -4646.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 554.4761999999999
$ws.Range("I122").Value = 335
$ws.Range("J122").Value = 993.4286
$ws.Range("K122").Value = 3015
$ws.Range("L122").Value = 8940.857399999999
$ws.Range("M122").Value = -565
$ws.Range("N122").Value = -13840.8574

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 990
$ws.Range("I125").Value = 990
$ws.Range("K125").Value = 2970
$ws.Range("M125").Value = 1950

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 14624.5
$ws.Range("I57").Value = 7166
$ws.Range("J57").Value = 37000
$ws.Range("K57").Value = 7166
$ws.Range("L57").Value = 37000
$ws.Range("M57").Value = -6346
$ws.Range("N57").Value = -38640

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3488.8333
$ws.Range("I126").Value = 3442.3635
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 10327.0905
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -7857.0905
$ws.Range("N126").Value = -16940

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5001631
$ws.Range("I132").Value = 5209991
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 15629973
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -15627443
$ws.Range("N132").Value = -8060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 15000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 15000
$ws.Range("K25").Value = 0
$ws.Range("L25").ClearContents()
$ws.Range("M25").Value = 15000
$ws.Range("N25").Value = -15460

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9999
$ws.Range("I40").Value = 9999
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 9999
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -9863

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 17321448
$ws.Range("I132").Value = 19726768
$ws.Range("J132").Value = 3139.8
$ws.Range("K132").Value = 59180304
$ws.Range("L132").Value = 9419.400000000001
$ws.Range("M132").Value = -59177774
$ws.Range("N132").Value = -14479.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1276.4445
$ws.Range("I126").Value = 1098.4
$ws.Range("J126").Value = 2166.6667
$ws.Range("K126").Value = 3295.2
$ws.Range("L126").Value = 6500.000100000001
$ws.Range("M126").Value = -825.2000000000003
$ws.Range("N126").Value = -11440.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11112569
$ws.Range("I132").Value = 12501462
$ws.Range("J132").Value = 1424.6
$ws.Range("K132").Value = 37504386
$ws.Range("L132").Value = 4273.799999999999
$ws.Range("M132").Value = -37501856
$ws.Range("N132").Value = -9333.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 14707823
$ws.Range("I136").Value = 15627037
$ws.Range("J136").Value = 400
$ws.Range("K136").Value = 46881111
$ws.Range("L136").Value = 1200
$ws.Range("M136").Value = -46878561
$ws.Range("N136").Value = -6300
